# Expand "Kế hoạch thực hiện đề tài" implementation-plan table from 7 to 14 rows.
$d = $word.ActiveDocument

# Locate the target table: the one whose header row starts with "STT".
$table = $null
for ($i = 1; $i -le $d.Tables.Count; $i++) {
    $candidate = $d.Tables.Item($i)
    if ($candidate.Cell(1, 1).Range.Text -like "STT*") {
        $table = $candidate
        break
    }
}

# --- Step 1 (row already present) ---------------------------------------
$table.Cell(2, 2).Range.Text = "Nghiên cứu công nghệ, phân tích đối thủ cạnh tranh"
$table.Cell(2, 3).Range.Text = "01/02 – 14/02/2026"
$table.Cell(2, 4).Range.Text = "BeeClass, Udemy, technology stack"

# --- Step 2 (row already present) ---------------------------------------
$table.Cell(3, 2).Range.Text = "Phân tích yêu cầu nghiệp vụ, use case modeling"
$table.Cell(3, 3).Range.Text = "15/02 – 21/02/2026"
$table.Cell(3, 4).Range.Text = "214 use cases, user stories"

# --- Step 3 (row already present) ---------------------------------------
$table.Cell(4, 2).Range.Text = "Thiết kế kiến trúc hệ thống (Hybrid Architecture)"
$table.Cell(4, 3).Range.Text = "22/02 – 28/02/2026"
$table.Cell(4, 4).Range.Text = "PlantUML diagrams, ADR"

# --- Step 4 (row already present) ---------------------------------------
$table.Cell(5, 2).Range.Text = "Thiết kế database schema & API specification"
$table.Cell(5, 3).Range.Text = "01/03 – 07/03/2026"
$table.Cell(5, 4).Range.Text = "ERD, Swagger/OpenAPI"

# --- Step 5 (row already present) ---------------------------------------
$table.Cell(6, 2).Range.Text = "Xây dựng Authentication & Authorization module"
$table.Cell(6, 3).Range.Text = "08/03 – 14/03/2026"
$table.Cell(6, 4).Range.Text = "JWT, RBAC, Spring Security"

# --- New steps 6-11: append rows right after the current last row (row 8 = old step 7) ---
$newSteps = @(
    @("6",  "Xây dựng Tenant Management & Billing System",       "15/03 – 21/03/2026", "Multi-tenant, VietQR payment"),
    @("7",  "Xây dựng Admin Dashboard & Auto-provisioning",      "22/03 – 28/03/2026", "Next.js, K8s API integration"),
    @("8",  "Xây dựng Course Service & Assignment Service",      "29/03 – 11/04/2026", "CRUD, business logic, DTOs"),
    @("9",  "Xây dựng Attendance Service & cross-service auth",  "12/04 – 18/04/2026", "Service-to-service JWT"),
    @("10", "Tích hợp AI Agent cho branding automation",         "19/04 – 25/04/2026", "GPT-4, DALL-E 3, cost `$0.19"),
    @("11", "Xây dựng Parent, Gamification, Forum Services",     "26/04 – 10/05/2026", "Unbundled pricing model")
)

# The original row 6 ("Kiểm thử toàn hệ thống...") is the anchor we must insert the
# new rows in front of, so it stays directly below the newly inserted block.
# `Rows.Add(before)` always inserts immediately above a fixed anchor, so walking
# the desired list back-to-front keeps the final on-page order correct.
$anchorRow = $table.Rows.Item(7)

for ($i = $newSteps.Length - 1; $i -ge 0; $i--) {
    $step = $newSteps[$i]
    $newRow = $table.Rows.Add($anchorRow)
    $newRow.Cells.Item(1).Range.Text = $step[0]
    $newRow.Cells.Item(2).Range.Text = $step[1]
    $newRow.Cells.Item(3).Range.Text = $step[2]
    $newRow.Cells.Item(4).Range.Text = $step[3]
}

# --- Former step 6 -> now step 12 (row index shifted by the 6 inserted rows) ---
$table.Cell(13, 1).Range.Text = "12"
$table.Cell(13, 2).Range.Text = "Unit testing (80% coverage) & Integration testing"
$table.Cell(13, 3).Range.Text = "11/05 – 17/05/2026"
$table.Cell(13, 4).Range.Text = "JUnit, Mockito, Testcontainers"

# --- Former step 7 -> now step 13 ---
$table.Cell(14, 1).Range.Text = "13"
$table.Cell(14, 2).Range.Text = "Load testing, performance tuning & deployment"
$table.Cell(14, 3).Range.Text = "18/05 – 25/05/2026"
$table.Cell(14, 4).Range.Text = "JMeter, AWS EKS production"

# --- New final step 14, appended at the end of the table ---
$lastRow = $table.Rows.Add()
$lastRow.Cells.Item(1).Range.Text = "14"
$lastRow.Cells.Item(2).Range.Text = "Hoàn thiện thesis report, slides, demo video"
$lastRow.Cells.Item(3).Range.Text = "26/05 – 31/05/2026"
$lastRow.Cells.Item(4).Range.Text = "Defense preparation"

Write-Output "Final row count: $($table.Rows.Count)"
